# Apply the "salesforce and brain tree paypal" edit:
# On the "GC_Codes" worksheet, five rows containing gift-certificate codes
# that are no longer valid/used are removed entirely (shifting the rows
# below them up), which also tidies up custom row heights and an
# unused cell style that only those rows used.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GC_Codes")

# Delete rows from the bottom up so row numbers of not-yet-deleted rows
# are not affected by earlier deletions.
$rowsToDelete = @(11, 10, 5, 4, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# The deleted rows were the only ones using a distinct (nearly-identical,
# unprotected) cell style; after their removal that style is unused.
# Normalize the remaining "code" rows (all except the bold header row A15)
# onto the style already used by A1:A2, so the redundant style can drop out.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A3:A14").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update selection to match the post-edit active cell recorded in the sheet.
$ws.Range("A7").Select()
